$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6: K6 gets a value of 0 (matching the B6:J6 "Total time minus" style)
# ---------------------------------------------------------------------------
$k6 = $ws.Range("K6")
$k6.HorizontalAlignment = -4152
$k6.Font.Size = 8
$k6.Font.Name = "Arial"
$k6.Font.Color = 0
$k6.Value = 0.0

# ---------------------------------------------------------------------------
# Row 7: label "10 k test" in A7 (left aligned, like the other row labels)
# and the new "10 k test" series values in B7:K7
# ---------------------------------------------------------------------------
$a7 = $ws.Range("A7")
$a7.HorizontalAlignment = -4131
$a7.Font.Size = 8
$a7.Font.Name = "Arial"
$a7.Font.Color = 0
$a7.Value = "10 k test"

$row7Values = @(0.9008, 0.8818, 0.895, 0.8944, 0.897, 0.8946, 0.8926, 0.8918, 0.8914, 0.891)
$row7Cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
for ($i = 0; $i -lt $row7Cols.Length; $i++) {
    $cell = $ws.Range($row7Cols[$i] + "7")
    $cell.HorizontalAlignment = -4152
    $cell.Font.Size = 8
    $cell.Font.Name = "Arial"
    $cell.Font.Color = 0
    $cell.Value = $row7Values[$i]
}

# ---------------------------------------------------------------------------
# Row 8: the "k" values (1 .. 10) used as the category axis for the new chart
# ---------------------------------------------------------------------------
$row8Values = @(1.0, 2.0, 3.0, 4.0, 5.0, 6.0, 7.0, 8.0, 9.0, 10.0)
$row8Cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
for ($i = 0; $i -lt $row8Cols.Length; $i++) {
    $cell = $ws.Range($row8Cols[$i] + "8")
    $cell.HorizontalAlignment = -4152
    $cell.Font.Size = 8
    $cell.Font.Name = "Arial"
    $cell.Font.Color = 0
    $cell.Value = $row8Values[$i]
}

# ---------------------------------------------------------------------------
# New chart #4: "K accuracy rate comparison" - plots B7:K7 against B8:K8
# ---------------------------------------------------------------------------
$chartObj4 = $ws.ChartObjects().Add(800, 800, 350, 200)
$chart4 = $chartObj4.Chart
$chart4.ChartType = 4
$chart4.HasTitle = $true
$chart4.ChartTitle.Text = "K accuracy rate comparison"
$series4 = $chart4.SeriesCollection().NewSeries()
$series4.Formula = "=SERIES('工作表1'!`$A`$7,'工作表1'!`$B`$8:`$L`$8,'工作表1'!`$B`$7:`$L`$7,1)"

# ---------------------------------------------------------------------------
# New chart #5: "Average Time" - plots B10:K10 against B1:K1
# ---------------------------------------------------------------------------
$chartObj5 = $ws.ChartObjects().Add(800, 1100, 350, 200)
$chart5 = $chartObj5.Chart
$chart5.ChartType = 4
$chart5.HasTitle = $true
$chart5.ChartTitle.Text = "Average Time"
$series5 = $chart5.SeriesCollection().NewSeries()
$series5.Formula = "=SERIES('工作表1'!`$A`$10,'工作表1'!`$B`$1:`$L`$1,'工作表1'!`$B`$10:`$L`$10,1)"
